# cell_values_and_formats.xlsx
# Commit: "add decimal place info for floats on reading from sheet"
#
# A new worksheet "zero_float" is inserted right after the "float" sheet.
# It mirrors the "float" sheet's layout/styles, but exercises a float
# value that happens to be a whole number (123, displayed as "123.00"
# under a 0.00 number format) -- i.e. a float with zero decimal places --
# to test decimal-place handling for such values.

$wb = $excel.ActiveWorkbook
$floatSheet = $wb.Worksheets.Item("float")

# Insert the new sheet right after "float". Worksheets.Add() makes the new
# sheet active/tabSelected and clears tabSelected on the previously active
# sheet, which is exactly what the source diff shows (sheet "string" loses
# tabSelected, "zero_float" gains it).
$ws = $wb.Worksheets.Add($null, $floatSheet)
$ws.Name = "zero_float"

# ---- layout: match the "float" sheet -------------------------------------
# (copy B1:C1 and A2:C7 separately so an inert/empty A1 cell record isn't
# materialized -- A1 has no content or special formatting on any sheet)
$floatSheet.Range("B1:C1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$floatSheet.Range("A2:C7").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Columns.Item(1).ColumnWidth = $floatSheet.Columns.Item(1).ColumnWidth
$ws.Columns.Item(2).ColumnWidth = $floatSheet.Columns.Item(2).ColumnWidth
$ws.Columns.Item(3).ColumnWidth = $floatSheet.Columns.Item(3).ColumnWidth

$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72

# ---- content ---------------------------------------------------------------
# Row 1 - column headers (same as every other sheet)
$ws.Range("B1").Value = "raw"
$ws.Range("C1").Value = "formula"

# Row 2 - "general" / converts value to non-float
$ws.Range("A2").Value = "general"
$ws.Range("B2").Value = "N/A - converts value to non-float"
$ws.Range("C2").Value = "N/A - converts value to non-float"

# Row 3 - "text" representation of the zero-decimal float/its formula.
# These two literal strings start with digits / "=" respectively; the
# second one must stay literal TEXT (not be parsed as a formula), so it is
# produced indirectly: compute the string with a helper formula, then
# paste-special only the *value* into the target cell (keeps it a shared
# string, does not disturb the text format already pasted above).
$ws.Range("A3").Value = "text"
$ws.Range("B3").Value = "123.00"

$ws.Range("Z1").Formula = '="="&"122.41+0.59"'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$ws.Range("Z1").Clear() | Out-Null

# Row 4 - "number": a float whose value has no fractional part
$ws.Range("A4").Value = "number"
$ws.Range("B4").Value = 123
$ws.Range("C4").Formula = "=122.41+0.59"

# Row 5 - "date"
$ws.Range("A5").Value = "date"
$ws.Range("B5").Value = 1000000000.12
$ws.Range("C5").Formula = "=999999999.41+1.04"

# Row 6 - "time"
$ws.Range("A6").Value = "time"
$ws.Range("B6").Value = 1000000000.12
$ws.Range("C6").Formula = "=999999999.41+0.04"

# Row 7 - "percentage" / converts value to percentage
$ws.Range("A7").Value = "percentage"
$ws.Range("B7").Value = "N/A - converts value to percentage"
$ws.Range("C7").Value = "N/A - converts value to percentage"

$ws.Range("B5").Select() | Out-Null
